# This script applies the data updates made in the "Add files via upload" commit
# to the AfDD_2025_Annex_Table_Tab32 "Tab32" worksheet:
#  1. Refreshed indicator values (columns C-G) for several countries / aggregates.
#  2. Corrected footnote text in the notes section (A103, A104): fixed mojibake
#     (UTF-8 double-encoding) characters and updated wording of the source note.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab32")

# --- Updated data values ---
$ws.Range("C6").Value = 89.141
$ws.Range("D6").Value = 89.9
$ws.Range("E6").Value = 88.548
$ws.Range("C10").Value = 34.901
$ws.Range("D10").Value = 34.648
$ws.Range("E10").Value = 35.107
$ws.Range("C13").Value = 74.333556
$ws.Range("D13").Value = 76.146556
$ws.Range("E13").Value = 72.531778
$ws.Range("F13").Value = 30.068998
$ws.Range("G13").Value = 35.884971
$ws.Range("F23").Value = 36.023464
$ws.Range("G23").Value = 40.278006
$ws.Range("C31").Value = 84.269
$ws.Range("D31").Value = 86.668
$ws.Range("E31").Value = 82.171
$ws.Range("C32").Value = 14.222
$ws.Range("D32").Value = 5.84
$ws.Range("E32").Value = 23.017
$ws.Range("C38").Value = 75.432417
$ws.Range("D38").Value = 76.871083
$ws.Range("E38").Value = 74.847
$ws.Range("F38").Value = 28.82903
$ws.Range("G38").Value = 37.791367
$ws.Range("C61").Value = 89.778923
$ws.Range("D61").Value = 92.056385
$ws.Range("E61").Value = 87.559231
$ws.Range("F61").Value = 48.963671
$ws.Range("G61").Value = 51.462883
$ws.Range("C62").Value = 80.322949
$ws.Range("D62").Value = 81.687949
$ws.Range("E62").Value = 78.998333
$ws.Range("F62").Value = 35.358739
$ws.Range("G62").Value = 39.739009
$ws.Range("C63").Value = 37.622138
$ws.Range("D63").Value = 36.903138
$ws.Range("E63").Value = 38.029915
$ws.Range("F63").Value = 16.282703
$ws.Range("G63").Value = 19.417129
$ws.Range("C64").Value = 52.89068
$ws.Range("D64").Value = 50.63776
$ws.Range("E64").Value = 54.6614
$ws.Range("F64").Value = 30.042083
$ws.Range("G64").Value = 36.293431
$ws.Range("C65").Value = 71.84305
$ws.Range("D65").Value = 71.1104
$ws.Range("E65").Value = 71.764
$ws.Range("F65").Value = 13.722841
$ws.Range("G65").Value = 18.053773
$ws.Range("C66").Value = 50.143429
$ws.Range("D66").Value = 50.035526
$ws.Range("E66").Value = 50.043211
$ws.Range("F66").Value = 17.368923
$ws.Range("G66").Value = 20.556808
$ws.Range("C67").Value = 73.047647
$ws.Range("D67").Value = 73.670647
$ws.Range("E67").Value = 72.421647
$ws.Range("F67").Value = 30.273788
$ws.Range("G67").Value = 34.448528
$ws.Range("C68").Value = 85.493789
$ws.Range("D68").Value = 87.206737
$ws.Range("E68").Value = 83.850789
$ws.Range("F68").Value = 39.042414
$ws.Range("G68").Value = 42.181919
$ws.Range("C69").Value = 91.261
$ws.Range("D69").Value = 93.6316
$ws.Range("E69").Value = 89.0566
$ws.Range("F69").Value = 33.713772
$ws.Range("G69").Value = 40.540259
$ws.Range("C70").Value = 92.535
$ws.Range("D70").Value = 94.72775
$ws.Range("E70").Value = 90.49175
$ws.Range("F70").Value = 37.326746
$ws.Range("G70").Value = 42.638319
$ws.Range("C71").Value = 89.778923
$ws.Range("D71").Value = 92.056385
$ws.Range("E71").Value = 87.559231
$ws.Range("F71").Value = 48.963671
$ws.Range("G71").Value = 51.462883
$ws.Range("F72").Value = 26.200295
$ws.Range("G72").Value = 33.886662
$ws.Range("C73").Value = 72.824333
$ws.Range("D73").Value = 73.414867
$ws.Range("E73").Value = 72.1976
$ws.Range("F73").Value = 32.056121
$ws.Range("G73").Value = 38.252636
$ws.Range("F75").Value = 38.28783
$ws.Range("G75").Value = 45.066483
$ws.Range("C76").Value = 71.643571
$ws.Range("D76").Value = 71.772286
$ws.Range("E76").Value = 71.655
$ws.Range("F76").Value = 22.494389
$ws.Range("G76").Value = 27.092184
$ws.Range("C77").Value = 54.801615
$ws.Range("D77").Value = 53.507154
$ws.Range("E77").Value = 55.623385
$ws.Range("F77").Value = 30.654612
$ws.Range("G77").Value = 37.664409
$ws.Range("C78").Value = 2.535593
$ws.Range("D78").Value = 2.857111
$ws.Range("E78").Value = 2.255667
$ws.Range("F78").Value = 18.267092
$ws.Range("G78").Value = 20.128951
$ws.Range("C79").Value = 9.487594
$ws.Range("D79").Value = 10.129875
$ws.Range("E79").Value = 9.034094
$ws.Range("F79").Value = 14.13205
$ws.Range("G79").Value = 15.780661
$ws.Range("F80").Value = 32.380161
$ws.Range("G80").Value = 37.672087
$ws.Range("C81").Value = 52.131875
$ws.Range("D81").Value = 47.87875
$ws.Range("E81").Value = 53.596875
$ws.Range("F81").Value = 28.886332
$ws.Range("G81").Value = 33.235701
$ws.Range("C82").Value = 79.723114
$ws.Range("D82").Value = 81.023143
$ws.Range("E82").Value = 78.427971
$ws.Range("F82").Value = 35.782934
$ws.Range("G82").Value = 40.066216
$ws.Range("C83").Value = 36.272395
$ws.Range("D83").Value = 35.882151
$ws.Range("E83").Value = 36.581826
$ws.Range("F83").Value = 15.222577
$ws.Range("G83").Value = 18.251353
$ws.Range("C84").Value = 92.718062
$ws.Range("D84").Value = 94.817
$ws.Range("E84").Value = 90.6625
$ws.Range("F84").Value = 33.204284
$ws.Range("G84").Value = 40.832019
$ws.Range("C86").Value = 78.876882
$ws.Range("D86").Value = 80.922529
$ws.Range("E86").Value = 76.953882
$ws.Range("F86").Value = 39.186847
$ws.Range("G86").Value = 42.468194
$ws.Range("C87").Value = 74.943889
$ws.Range("D87").Value = 74.6905
$ws.Range("E87").Value = 74.855278
$ws.Range("F87").Value = 19.675096
$ws.Range("G87").Value = 25.924583
$ws.Range("C88").Value = 52.19125
$ws.Range("D88").Value = 49.94725
$ws.Range("E88").Value = 53.88
$ws.Range("F88").Value = 25.241447
$ws.Range("G88").Value = 30.948869
$ws.Range("C89").Value = 49.622613
$ws.Range("D89").Value = 47.557387
$ws.Range("E89").Value = 50.778516
$ws.Range("F89").Value = 16.677085
$ws.Range("G89").Value = 20.855774
$ws.Range("C90").Value = 11.794364
$ws.Range("D90").Value = 11.468364
$ws.Range("E90").Value = 12.101864
$ws.Range("F90").Value = 15.195559
$ws.Range("G90").Value = 16.887058
$ws.Range("C91").Value = 90.291269
$ws.Range("D91").Value = 92.9755
$ws.Range("E91").Value = 88.131923
$ws.Range("F91").Value = 33.789683
$ws.Range("G91").Value = 41.999554
$ws.Range("F92").Value = 25.808621
$ws.Range("G92").Value = 37.366015
$ws.Range("C93").Value = 59.1064
$ws.Range("D93").Value = 57.3484
$ws.Range("E93").Value = 60.719
$ws.Range("F93").Value = 24.487261
$ws.Range("G93").Value = 26.696562
$ws.Range("C94").Value = 47.785842
$ws.Range("D94").Value = 44.054421
$ws.Range("E94").Value = 50.597842
$ws.Range("F94").Value = 16.587724
$ws.Range("G94").Value = 21.026295
$ws.Range("C95").Value = 86.646083
$ws.Range("D95").Value = 88.456333
$ws.Range("E95").Value = 85.07425
$ws.Range("F95").Value = 32.446778
$ws.Range("G95").Value = 40.105879
$ws.Range("C96").Value = 62.347273
$ws.Range("D96").Value = 62.116
$ws.Range("E96").Value = 62.626909
$ws.Range("C97").Value = 88.406296
$ws.Range("D97").Value = 91.305593
$ws.Range("E97").Value = 86.029519
$ws.Range("F97").Value = 42.209394
$ws.Range("G97").Value = 47.28168
$ws.Range("F98").Value = 25.709208
$ws.Range("G98").Value = 31.790925

# --- Corrected footnote / source text ---
$ws.Range("A103").Value = 'Regional Economic Communities: CEN-SAD = "Community of Sahel-Saharan States"; COMESA = "Common Market for Eastern and Southern Africa"; EAC = "East African Community"; ECCAS = "Economic Community of Central African States"; ECOWAS = "Economic Community of West African States"; IGAD = "Intergovernmental Authority on Development"; SADC = "Southern African Development Community"; UMA = "Arab Maghreb Union"; PALOP = "Países Africanos de Língua Oficial Portuguesa"; ASEAN = "Association of Southeast Asian Nations"; MERCOSUR = "Mercado Común del Sur". EU27 = "European Union (27 members)". OECD = "Organisation for Economic Co-operation and Development".'
$ws.Range("A104").Value = 'Source: International Labour Organization - ILOSTAT (retrieved 09/09/2025), "The Informality Database" (Elgin, C., M. A. Kose, F. Ohnsorge, and S. Yu. 2021. Understanding Informality. CERP Discussion Paper 16497, Centre for Economic Policy Research, London - data updated 09/01/2024).'
